$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new (empty) column before the existing column D.
# This pushes the old "D" formula column (CONCATENATE results) to column E,
# while columns A, B and C (and their content/styles) stay where they are.
$ws.Columns("D").Insert()

# Type the new header value "Code" into C3 (below the "Familly" header in B3),
# copying the same header formatting (bold font, fill, border) used by B3.
$ws.Range("B3").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C3").Value = "Code"

# Restore the view: scroll back to the top-left and select a single cell (C5),
# matching the saved view state after this edit.
$ws.Range("C5").Select()
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
